$wb = $excel.ActiveWorkbook

$updates = @{
    "F3"  = 299
    "F5"  = 23
    "F6"  = 316
    "F7"  = 9440
    "F8"  = 80
    "F10" = 133
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
